# Implantação da busca pela home com sucesso através de massa de dados
#
# - Renomeia a aba "buscarHome" para "buscarHomeSucesso" e remove a coluna
#   de resultado (B1) que não é mais utilizada.
# - Corrige o texto do modelo de "BOSE SOUNDLINK BLUETOOTH SPEAKER III"
#   para "Bose Soundlink Bluetooth Speaker III".
# - Acrescenta uma nova aba "buscarHomeFail" para a massa de dados de busca
#   com falha, ficando como a aba ativa.

$wb = $excel.ActiveWorkbook

# --- buscarHome -> buscarHomeSucesso -------------------------------------
$wsSucesso = $wb.Worksheets.Item("buscarHome")
$wsSucesso.Name = "buscarHomeSucesso"

# Corrige a capitalização do nome do produto (mesma shared string usada em A2)
$wsSucesso.Range("A2").Value = "Bose Soundlink Bluetooth Speaker III"

# A coluna "Resultado" (B1) não é mais necessária nesta massa de dados
$wsSucesso.Range("B1").ClearContents()

# --- nova aba buscarHomeFail ----------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFail = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsFail.Name = "buscarHomeFail"

$psFail = $wsFail.PageSetup
$psFail.LeftMargin = 0.511811024 * 72
$psFail.RightMargin = 0.511811024 * 72
$psFail.TopMargin = 0.78740157499999996 * 72
$psFail.BottomMargin = 0.78740157499999996 * 72
$psFail.HeaderMargin = 0.31496062000000002 * 72
$psFail.FooterMargin = 0.31496062000000002 * 72

# Seleção final: B1 na aba de sucesso (célula vazia após remover "Resultado")
$wsSucesso.Range("B1").Select() | Out-Null

# A aba de falha fica selecionada/ativa por último
$wsFail.Activate() | Out-Null

Write-Host "buscarHome -> buscarHomeSucesso + buscarHomeFail aplicado"
